$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = "[1, 0, 0, 1, 0, 0, 0]"
$ws.Range("E3").Value = "['Normal', 'ParamViolation']"

# Row 12
$ws.Range("D12").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E12").Value = "['Normal', 'HardwareFault']"

# Row 53
$ws.Range("D53").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E53").Value = "['Normal', 'HardwareFault']"

# Row 73
$ws.Range("D73").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E73").Value = "['Normal']"

# Row 82
$ws.Range("D82").Value = "[1, 1, 1, 0, 0, 0, 0]"
$ws.Range("E82").Value = "['Normal', 'SurroundingEnvironment', 'HardwareFault']"

# Row 92
$ws.Range("D92").Value = "[1, 0, 1, 0, 0, 0, 1]"
$ws.Range("E92").Value = "['Normal', 'HardwareFault', 'SoftwareFault']"

# Row 109
$ws.Range("D109").Value = "[1, 1, 0, 0, 0, 0, 1]"
$ws.Range("E109").Value = "['Normal', 'SurroundingEnvironment', 'SoftwareFault']"
